# "Generate Report for Archive"
#
# The localization status report is regenerated: the three tracked files
# db1cacb9-*, e3a9b16a-* and fd0329ee-* get reordered (e3a9b16a now comes
# first, then fd0329ee, then db1cacb9) and the status of e3a9b16a-* /
# fd0329ee-* moves from "Ready for handoff" to "In Translation" (with
# updated handoff timestamps). Hyperlink addresses (r:id -> target URL)
# stay anchored to their original cell position; only the visible display
# text is refreshed to track the new cell contents.

function Set-HyperlinkDisplay($ws, $ref, $newText) {
    foreach ($hl in $ws.Hyperlinks) {
        $r = $hl.Range.Address()
        if ($r -eq $ref) {
            $hl.TextToDisplay = $newText
        }
    }
}

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Sheet "Overview"
# ---------------------------------------------------------------------
$ws = $wb.Worksheets.Item("Overview")

$ws.Range("A3").Value = "e3a9b16a-1b58-4003-9522-9be81fb79950.md"
$ws.Range("B3").Value = "In Translation"
$ws.Range("C3").Value = "In Translation"
$ws.Range("D3").Value = "2016-03-23 18:42:44"

$ws.Range("A4").Value = "fd0329ee-f51a-4870-bb82-0f99a9481b66.md"
$ws.Range("B4").Value = "In Translation"
$ws.Range("C4").Value = "In Translation"
$ws.Range("D4").Value = "2016-03-23 18:42:44"

$ws.Range("A5").Value = "db1cacb9-7b90-49aa-8fbc-172b5f18c289.md"
$ws.Range("B5").Value = "Ready for handoff"
$ws.Range("C5").Value = "Ready for handoff"
$ws.Range("D5").Value = "2016-03-23 18:40:28"

Set-HyperlinkDisplay $ws '$A$3' "e3a9b16a-1b58-4003-9522-9be81fb79950.md"
Set-HyperlinkDisplay $ws '$A$4' "fd0329ee-f51a-4870-bb82-0f99a9481b66.md"
Set-HyperlinkDisplay $ws '$A$5' "db1cacb9-7b90-49aa-8fbc-172b5f18c289.md"

# ---------------------------------------------------------------------
# Sheet "zh-cn"
# ---------------------------------------------------------------------
$ws = $wb.Worksheets.Item("zh-cn")

$ws.Range("A3").Value = "e3a9b16a-1b58-4003-9522-9be81fb79950.md"
$ws.Range("C3").Value = "In Translation"
$ws.Range("D3").Value = "e3a9b16a-1b58-4003-9522-9be81fb79950.8ea1b1c973f2d86da1c4c742ae52735de00dfee0.zh-cn.xlf"
$ws.Range("E3").Value = "2016-03-23 18:42:40"

$ws.Range("A4").Value = "fd0329ee-f51a-4870-bb82-0f99a9481b66.md"
$ws.Range("C4").Value = "In Translation"
$ws.Range("D4").Value = "fd0329ee-f51a-4870-bb82-0f99a9481b66.1da723cd92c9571137e1079134d51d53b326f9ce.zh-cn.xlf"
$ws.Range("E4").Value = "2016-03-23 18:42:40"

$ws.Range("A5").Value = "db1cacb9-7b90-49aa-8fbc-172b5f18c289.md"
$ws.Range("C5").Value = "Ready for handoff"
$ws.Range("D5").Value = "db1cacb9-7b90-49aa-8fbc-172b5f18c289.53e5b2f026856f39bf947f8f1678af1ee414f37d.zh-cn.xlf"
$ws.Range("E5").Value = "2016-03-23 18:40:24"

Set-HyperlinkDisplay $ws '$A$3' "e3a9b16a-1b58-4003-9522-9be81fb79950.md"
Set-HyperlinkDisplay $ws '$D$3' "e3a9b16a-1b58-4003-9522-9be81fb79950.8ea1b1c973f2d86da1c4c742ae52735de00dfee0.zh-cn.xlf"
Set-HyperlinkDisplay $ws '$A$4' "fd0329ee-f51a-4870-bb82-0f99a9481b66.md"
Set-HyperlinkDisplay $ws '$D$4' "fd0329ee-f51a-4870-bb82-0f99a9481b66.1da723cd92c9571137e1079134d51d53b326f9ce.zh-cn.xlf"
Set-HyperlinkDisplay $ws '$A$5' "db1cacb9-7b90-49aa-8fbc-172b5f18c289.md"
Set-HyperlinkDisplay $ws '$D$5' "db1cacb9-7b90-49aa-8fbc-172b5f18c289.53e5b2f026856f39bf947f8f1678af1ee414f37d.zh-cn.xlf"

# ---------------------------------------------------------------------
# Sheet "de-de"
# ---------------------------------------------------------------------
$ws = $wb.Worksheets.Item("de-de")

$ws.Range("A3").Value = "e3a9b16a-1b58-4003-9522-9be81fb79950.md"
$ws.Range("C3").Value = "In Translation"
$ws.Range("D3").Value = "e3a9b16a-1b58-4003-9522-9be81fb79950.8ea1b1c973f2d86da1c4c742ae52735de00dfee0.de-de.xlf"
$ws.Range("E3").Value = "2016-03-23 18:42:44"

$ws.Range("A4").Value = "fd0329ee-f51a-4870-bb82-0f99a9481b66.md"
$ws.Range("C4").Value = "In Translation"
$ws.Range("D4").Value = "fd0329ee-f51a-4870-bb82-0f99a9481b66.1da723cd92c9571137e1079134d51d53b326f9ce.de-de.xlf"
$ws.Range("E4").Value = "2016-03-23 18:42:44"

$ws.Range("A5").Value = "db1cacb9-7b90-49aa-8fbc-172b5f18c289.md"
$ws.Range("C5").Value = "Ready for handoff"
$ws.Range("D5").Value = "db1cacb9-7b90-49aa-8fbc-172b5f18c289.53e5b2f026856f39bf947f8f1678af1ee414f37d.de-de.xlf"
$ws.Range("E5").Value = "2016-03-23 18:40:28"

Set-HyperlinkDisplay $ws '$A$3' "e3a9b16a-1b58-4003-9522-9be81fb79950.md"
Set-HyperlinkDisplay $ws '$D$3' "e3a9b16a-1b58-4003-9522-9be81fb79950.8ea1b1c973f2d86da1c4c742ae52735de00dfee0.de-de.xlf"
Set-HyperlinkDisplay $ws '$A$4' "fd0329ee-f51a-4870-bb82-0f99a9481b66.md"
Set-HyperlinkDisplay $ws '$D$4' "fd0329ee-f51a-4870-bb82-0f99a9481b66.1da723cd92c9571137e1079134d51d53b326f9ce.de-de.xlf"
Set-HyperlinkDisplay $ws '$A$5' "db1cacb9-7b90-49aa-8fbc-172b5f18c289.md"
Set-HyperlinkDisplay $ws '$D$5' "db1cacb9-7b90-49aa-8fbc-172b5f18c289.53e5b2f026856f39bf947f8f1678af1ee414f37d.de-de.xlf"
